$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.773.81'
$ws.Range("E2").Value = '  +1.59%  '

$ws.Range("D3").Value = '3.738.08'
$ws.Range("E3").Value = '  -1.76%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.30'
$ws.Range("E5").Value = '  +0.73%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.32'
$ws.Range("E6").Value = '  -4.94%  '

$ws.Range("D7").Value = '3.736.68'
$ws.Range("E7").Value = '  -1.77%  '

$ws.Range("E9").Value = '  +1.66%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.164'
$ws.Range("E10").Value = '  +2.82%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.37'
$ws.Range("E11").Value = '  +2.76%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.459'
$ws.Range("E12").Value = '  -0.97%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.01'
$ws.Range("E13").Value = '  -0.75%  '

$ws.Range("E14").Value = '  -0.48%  '

$ws.Range("D15").Value = '4.363.18'
$ws.Range("E15").Value = '  -1.75%  '

$ws.Range("D16").Value = '3.735.94'
$ws.Range("E16").Value = '  -1.78%  '

$ws.Range("D17").Value = '68.761.00'
$ws.Range("E17").Value = '  +1.59%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.24'
$ws.Range("E18").Value = '  +0.55%  '

$ws.Range("E19").Value = '  +0.33%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.26'
$ws.Range("E20").Value = '  +3.89%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '496.65'
$ws.Range("E21").Value = '  +1.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.04'
$ws.Range("E22").Value = '  +10.43%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.723'
$ws.Range("E23").Value = '  -2.89%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.87'
$ws.Range("E24").Value = '  -0.40%  '

$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000142'
$ws.Range("E25").Value = '  -6.13%  '

$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.31'
$ws.Range("E26").Value = '  -2.82%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.38'
$ws.Range("E27").Value = '  +0.71%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.07'
$ws.Range("E28").Value = '  -1.24%  '

$ws.Range("E29").Value = '  +0.03%  '

$ws.Range("E30").Value = '  -0.68%  '

$ws.Range("E31").Value = '  -0.06%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.91'
$ws.Range("E32").Value = '  +3.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.69'
$ws.Range("E33").Value = '  -2.20%  '

$ws.Range("D34").Value = '3.883.38'
$ws.Range("E34").Value = '  -1.64%  '

$ws.Range("B35").Value = 'RenzoRestakedETH'
$ws.Range("C35").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D35").Value = '3.670.52'
$ws.Range("E35").Value = '  -1.91%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.108'
$ws.Range("E36").Value = '  -0.60%  '

$ws.Range("E37").Value = '  +0.02%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.02'
$ws.Range("E38").Value = '  +0.91%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.82'
$ws.Range("E39").Value = '  +0.36%  '

$ws.Range("E40").Value = '  -2.06%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.325'
$ws.Range("E41").Value = '  -1.74%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '433.03'
$ws.Range("E42").Value = '  -4.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '49.18'
$ws.Range("E43").Value = '  +0.02%  '

$ws.Range("E44").Value = '  -1.55%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.87'
$ws.Range("E45").Value = '  -0.41%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.42'
$ws.Range("E46").Value = '  +0.89%  '

$ws.Range("E47").Value = '  +0.00%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.65'
$ws.Range("E48").Value = '  -1.73%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.09'
$ws.Range("E49").Value = '  +0.88%  '

$ws.Range("E50").Value = '  +0.93%  '

$ws.Range("D51").Value = '2.745.29'
$ws.Range("E51").Value = '  -3.51%  '
